# PPE-105234 : Modified tests to handle qa02 as perf (Sub Tasks : PPE-118033, PPE-118034)
#
# Row 3 (A3/B3) already held the "testingpurpose/2019/6" bulk-import URLs with
# hyperlinks. This adds a second "testingpurpose" (no date) variant in row 2,
# turning A2/B2 into live hyperlinks with the same Hyperlink cell style as
# row 3, and moves the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fitnessTestingPurpose = "http://www.qa00.webmd.com/fitness-exercise/features/tired-of-exercise/testingpurpose"
$menTestingPurpose      = "http://www.qa00.webmd.com/men/stick-with-fitness-plan/testingpurpose"

# Update the row 2 cell text to the new "testingpurpose" (no date) URLs.
$ws.Range("A2").Value = $fitnessTestingPurpose
$ws.Range("B2").Value = $menTestingPurpose

# Turn them into real hyperlinks (targets match the new cell text).
$ws.Hyperlinks.Add($ws.Range("A2"), $fitnessTestingPurpose)
$ws.Hyperlinks.Add($ws.Range("B2"), $menTestingPurpose)

# Hyperlinks.Add stamps its own style variant; re-apply the shared
# "Hyperlink" cell style so A2/B2 match A3/B3 exactly.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"

# Move the selection cursor to B3 (was A4).
$ws.Range("B3").Select()
